$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, shifting existing rows 12-19 down to 13-20
$ws.Rows("12:12").Insert()

# Populate the new row 12 with the new weekly record
$ws.Range("A12").Value2 = 1
$ws.Range("B12").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C12").Value2 = "Arica y Parinacota"
$ws.Range("D12").Value2 = 45118
$ws.Range("D12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E12").Value2 = 15
$ws.Range("F12").Value2 = 100112017
$ws.Range("G12").Value2 = "Ramas de apio"
$ws.Range("H12").Value2 = "Sin especificar"
$ws.Range("I12").Value2 = "Primera"
$ws.Range("J12").Value2 = 200
$ws.Range("K12").Value2 = 4000
$ws.Range("L12").Value2 = 5000
$ws.Range("M12").Value2 = 4500
$ws.Range("N12").Value2 = "$/atado 7 kilos"
$ws.Range("O12").Value2 = "Región de Arica y Parinacota"
$ws.Range("P12").Value2 = 4500
$ws.Range("Q12").Value2 = 1
$ws.Range("R12").Value2 = "Hortaliza"
